$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the regression output values to their rounded (2 decimal) forms
$ws.Range("B2").Value = "-0.37***"
$ws.Range("B3").Value = "-3.46***"
$ws.Range("C3").Value = "-0.81***"
